$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4061793025036406
$ws.Range("C2").Value = 0.9919907002636758
$ws.Range("D2").Value = 0.5025197286154568
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(learning_rate=0.5))])"
$ws.Range("G2").Value = 0.122648122766744
$ws.Range("H2").Value = 0.991
